$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(168).Insert()
$ws.Rows.Item(168).Insert()

$ws.Range("A168").Value = 1
$ws.Range("B168").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C168").Value = "Arica y Parinacota"
$ws.Range("D168").Value = 45119
$ws.Range("E168").Value = 15
$ws.Range("F168").Value = "Fruta"
$ws.Range("G168").Value = 100106
$ws.Range("H168").Value = "Oleaginosos"
$ws.Range("I168").Value = 100106002
$ws.Range("J168").Value = "Palta"
$ws.Range("K168").Value = "Hass"
$ws.Range("L168").Value = "Primera"
$ws.Range("M168").Value = 208
$ws.Range("N168").Value = 24000
$ws.Range("O168").Value = 25000
$ws.Range("P168").Value = 24500
$ws.Range("Q168").Value = "$/bandeja 10 kilos"
$ws.Range("R168").Value = "Perú"
$ws.Range("S168").Value = 2450
$ws.Range("T168").Value = 10

$ws.Range("A169").Value = 1
$ws.Range("B169").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C169").Value = "Arica y Parinacota"
$ws.Range("D169").Value = 45119
$ws.Range("E169").Value = 15
$ws.Range("F169").Value = "Fruta"
$ws.Range("G169").Value = 100106
$ws.Range("H169").Value = "Oleaginosos"
$ws.Range("I169").Value = 100106002
$ws.Range("J169").Value = "Palta"
$ws.Range("K169").Value = "Hass"
$ws.Range("L169").Value = "Segunda"
$ws.Range("M169").Value = 208
$ws.Range("N169").Value = 21000
$ws.Range("O169").Value = 22000
$ws.Range("P169").Value = 21500
$ws.Range("Q169").Value = "$/bandeja 10 kilos"
$ws.Range("R169").Value = "Perú"
$ws.Range("S169").Value = 2150
$ws.Range("T169").Value = 10
